$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.431.13"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.572.16"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3729"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.94"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3391"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07572"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.139"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.33"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.990"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.965"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "1.575.53"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.90"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.33"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").Value = "22.424.72"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.652"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.25"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.010"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.45"
$ws.Range("D31").Value = "1.752.71"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.045"
$ws.Range("E32").Value = "  +4.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.154"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.795"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08392"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02476"
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2281"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06516"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.461"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6221"
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.92"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.813"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5791"
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.46"
$ws.Range("E48").Value = "  +3.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.073"
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("E50").Value = "  -6.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07323"
$ws.Range("E51").Value = "  -0.14%  "
